# Apply the "add examples to FamilyMemberHistory, and cleanup narrative" edit:
#  - fix casing of the extension name/URL (ethnicity -> Ethnicity)
#  - bump the StructureDefinition Date
#  - fix casing in the narrative Description / Definition text
#  - point the binding Value Set at the real v3-Ethnicity value set
#  - narrow the "Binding Value Set" column on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- URL (Metadata!B2, also used verbatim as Extension.url's Fixed Value
#     on the Elements sheet, Q5) ---
$newUrl = "https://johnmoehrke.github.io/MHV-PGHD/StructureDefinition/FM-Ethnicity"
$meta.Range("B2").Value = $newUrl
$elements.Range("Q5").Value = $newUrl

# --- Name / Title (Metadata!B4 & B5) and the Extension row's Short
#     description (Elements!K2) all shared the same text "ethnicity" ---
$newName = "Ethnicity"
$meta.Range("B4").Value = $newName
$meta.Range("B5").Value = $newName
$elements.Range("K2").Value = $newName

# --- Date (Metadata!B8) ---
$meta.Range("B8").Value = "2022-04-11T07:37:02-05:00"

# --- Description (Metadata!B12) / Definition (Elements!L2) ---
$newDescription = "What is the Ethnicity of this family member`n`nNote would like to use the us-core defined extension, but it is not allowed in FamilyMemberHistory [jira issue](https://jira.hl7.org/browse/FHIR-35997)"
$meta.Range("B12").Value = $newDescription
$elements.Range("L2").Value = $newDescription

# --- Binding Value Set (Elements!Y7) ---
$elements.Range("Y7").Value = "http://terminology.hl7.org/ValueSet/v3-Ethnicity"

# --- Narrow column Y ("Binding Value Set") on the Elements sheet from
#     ~55.6 down to ~45 characters wide ---
$elements.Columns.Item(25).ColumnWidth = 44.14
